$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the PEPMatch benchmarking numbers (row 2) ---
$ws.Range("B2").Value2 = 46.247
$ws.Range("D2").Value2 = 193.638
$ws.Range("E2").Value2 = 239.885

# --- Reset the formatting on B2:F2 back to the plain/default style ---
$clear = $ws.Range("B2:F2")
$clear.Font.Name = "Calibri"
$clear.Font.Size = 11
$clear.Borders.LineStyle = -4142
$clear.HorizontalAlignment = 1
$clear.VerticalAlignment = -4107
$clear.NumberFormat = "General"

# --- Update the selection to match the new edit location ---
$ws.Range("B2:F2").Select() | Out-Null
